$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 782 (pushes existing rows 782-823 down to 783-824,
# and Excel auto-extends the used range / dimension to A1:D824).
$ws.Rows.Item(782).Insert()

# New row 782 data: 2026/02/12, 木 (Thursday), 2, 47
# The date-like text must be written as literal text (matches the other date
# cells in the sheet, which are inline/shared strings, not real dates). Using
# Range.Value directly on a date-shaped string causes Excel to auto-convert it
# into a date serial number + date number format, so instead we write it as a
# text formula result and then paste-special just the value back over itself,
# which keeps it as plain text without adding any new cell style.
$ws.Range("A782").Formula = "=""2026/02/12"""
$ws.Range("A782").Copy()
$ws.Range("A782").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B782").Value = "木"
$ws.Range("C782").Value = 2
$ws.Range("D782").Value = 47
